# Japanese (ja-JP) translation workbook:
# Add a new translation key "strChkAbsoluteIntegral" for the string
# "Compute the absolute-value integral?" right before the existing
# "strChkComputeDerivative" entry (i.e. as the 7th data row of the
# "Tabla13" table, which lives on the sole worksheet "ja-JP").
# This pushes every following row down by one and grows the table /
# worksheet extent from B2:E167 to B2:E168.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ja-JP")
$lo = $ws.ListObjects.Item("Tabla13")

# Insert a brand-new row inside the table's data body, right before the
# row that currently holds "strChkComputeDerivative" (the 7th data row,
# i.e. worksheet row 9). Inserting through a row that belongs to the
# table body keeps the formatting consistent with its neighbours and
# shifts every subsequent row (and the table's bounds) down by one.
$bodyRange = $lo.DataBodyRange
$targetRow = $bodyRange.Rows.Item(7)
$targetRow.Insert()

# The worksheet's own dimension grows automatically with the inserted
# row, but the table definition needs to be told explicitly to include
# it too, so grow the table by one row while keeping the same width.
$oldRange = $lo.Range
$newRange = $oldRange.Resize($oldRange.Rows.Count + 1, $oldRange.Columns.Count)
$lo.Resize($newRange)

# Populate the new row (Key / Comment / English columns). The ja.JP
# column (E) is intentionally left blank, same as the other untranslated
# rows in this sheet.
$ws.Range("B9").Value = "strChkAbsoluteIntegral"
$ws.Range("C9").Value = "In ""settings"" form, tab ""Integration"""
$ws.Range("D9").Value = "Compute the absolute-value integral?"
